# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have been generated:
#   - The "Status" text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears (Overview
#     per-locale columns, and the Status column on each locale sheet).
#   - Each locale sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns are populated for both data rows.
#   - New hyperlinks are added for the newly populated "Latest Target File"
#     cells.
#   - The columns that now hold longer text are widened to fit.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40c58f6a41fd41579fbd881ccb549029b88ef7a8/e2e/a.md"

$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$zhHandbackDateTime = "2016-09-04 12:40:50"
$deHandbackDateTime = "2016-09-04 12:40:57"

# ---------------------------------------------------------------------
# Overview sheet: update the per-locale status text (columns E and F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet: status text, handback columns, hyperlinks, widths
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("K2").Value = $zhHandbackDateTime

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = $zhHandbackFile
$wsZh.Range("K3").Value = $zhHandbackDateTime

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet: status text, handback columns, hyperlinks, widths
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("K2").Value = $deHandbackDateTime

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = $deHandbackFile
$wsDe.Range("K3").Value = $deHandbackDateTime

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1
